# Reposition Weekly/Quarterly report blocks on the Dashboard sheet so they
# sit below the six existing charts (and the alerts/recommendations
# sections), clearing the rows they used to occupy.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dashboard")

# ---------------------------------------------------------------------
# 1) WEEKLY REPORT  (old: A25:F34)  ->  (new: A88:F96)
# ---------------------------------------------------------------------

# -- copy the cell formatting to the new location -------------------------
$ws.Range("A25").Copy()
$ws.Range("A88").PasteSpecial($xlPasteFormats)

$ws.Range("A26").Copy()
$ws.Range("A89").PasteSpecial($xlPasteFormats)

$ws.Range("A28:F28").Copy()
$ws.Range("A90:F90").PasteSpecial($xlPasteFormats)

$ws.Range("A29:F33").Copy()
$ws.Range("A91:F95").PasteSpecial($xlPasteFormats)

$ws.Range("A34:C34").Copy()
$ws.Range("A96:C96").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

# -- write the relocated content/labels/formulas --------------------------
$ws.Range("A88").Value = "WEEKLY REPORT"
$ws.Range("A89").Value = "Week-by-week breakdown for selected month"

$ws.Range("A90").Value = "Week"
$ws.Range("B90").Value = "Incidents"
$ws.Range("C90").Value = "Ambulance"
$ws.Range("D90").Value = "Avg Age"
$ws.Range("E90").Value = "Top Injury"
$ws.Range("F90").Value = "Top Venue"

$weekRows = @(91, 92, 93, 94, 95)
for ($i = 0; $i -lt 5; $i++) {
    $r = $weekRows[$i]
    $wk = $i + 1
    $ws.Range("A$r").Value = "Week $wk"
    $ws.Range("B$r").Formula = "=SUMPRODUCT((MedicalData[Month]=`$B`$2)*((WEEKNUM(MedicalData[Date])-WEEKNUM(DATE(YEAR(MedicalData[Date]),MONTH(MedicalData[Date]),1))+1)=$wk)*1)"
    $ws.Range("C$r").Formula = "=SUMPRODUCT((MedicalData[Month]=`$B`$2)*((WEEKNUM(MedicalData[Date])-WEEKNUM(DATE(YEAR(MedicalData[Date]),MONTH(MedicalData[Date]),1))+1)=$wk)*(MedicalData[Hospital Transportation]=""Yes"")*1)"
    $ws.Range("D$r").Formula = "=IFERROR(ROUND(AGGREGATE(1,6,MedicalData[Age]/((MedicalData[Month]=`$B`$2)*((WEEKNUM(MedicalData[Date])-WEEKNUM(DATE(YEAR(MedicalData[Date]),MONTH(MedicalData[Date]),1))+1)=$wk))),0),""-"")"
    $ws.Range("E$r").Formula = "=IF(B$r=0,""-"",IFERROR(INDEX(MedicalData[Specific injuries treated],MATCH(1,(MedicalData[Month]=`$B`$2)*((WEEKNUM(MedicalData[Date])-WEEKNUM(DATE(YEAR(MedicalData[Date]),MONTH(MedicalData[Date]),1))+1)=$wk),0)),""-""))"
    $ws.Range("F$r").Formula = "=IF(B$r=0,""-"",IFERROR(INDEX(MedicalData[Site],MATCH(1,(MedicalData[Month]=`$B`$2)*((WEEKNUM(MedicalData[Date])-WEEKNUM(DATE(YEAR(MedicalData[Date]),MONTH(MedicalData[Date]),1))+1)=$wk),0)),""-""))"
}

$ws.Range("A96").Value = "TOTAL"
$ws.Range("B96").Formula = "=SUM(B91:B95)"
$ws.Range("C96").Formula = "=SUM(C91:C95)"
$ws.Range("D96").Formula = "=IFERROR(ROUND(AVERAGEIF(MedicalData[Month],`$B`$2,MedicalData[Age]),0),""-"")"

# ---------------------------------------------------------------------
# 2) QUARTERLY REPORT  (old: A37:F46)  ->  (new: A100:F109)
# ---------------------------------------------------------------------

$ws.Range("A37").Copy()
$ws.Range("A100").PasteSpecial($xlPasteFormats)

$ws.Range("B38").Copy()
$ws.Range("B101").PasteSpecial($xlPasteFormats)

$ws.Range("A40:F40").Copy()
$ws.Range("A103:F103").PasteSpecial($xlPasteFormats)

$ws.Range("F41").Copy()
$ws.Range("F104").PasteSpecial($xlPasteFormats)

$ws.Range("F42").Copy()
$ws.Range("F105").PasteSpecial($xlPasteFormats)

$ws.Range("F44").Copy()
$ws.Range("F107").PasteSpecial($xlPasteFormats)

$ws.Range("A46").Copy()
$ws.Range("A109").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false

$ws.Range("A100").Value = "QUARTERLY REPORT"

$ws.Range("A101").Value = "Quarter:"
$ws.Range("B101").Formula = '=IF($B$2="All Time","All",IF(OR(LEFT($B$2,3)="Jan",LEFT($B$2,3)="Feb",LEFT($B$2,3)="Mar"),"Q1",IF(OR(LEFT($B$2,3)="Apr",LEFT($B$2,3)="May",LEFT($B$2,3)="Jun"),"Q2",IF(OR(LEFT($B$2,3)="Jul",LEFT($B$2,3)="Aug",LEFT($B$2,3)="Sep"),"Q3","Q4"))))'

$ws.Range("A103").Value = "Metric"
$ws.Range("B103").Value = "Q1 (Jan-Mar)"
$ws.Range("C103").Value = "Q2 (Apr-Jun)"
$ws.Range("D103").Value = "Q3 (Jul-Sep)"
$ws.Range("E103").Value = "Q4 (Oct-Dec)"
$ws.Range("F103").Value = "YTD Total"

$ws.Range("A104").Value = "Total Incidents"
$ws.Range("B104").Formula = '=COUNTIF(MedicalData[Quarter],"Q1")'
$ws.Range("C104").Formula = '=COUNTIF(MedicalData[Quarter],"Q2")'
$ws.Range("D104").Formula = '=COUNTIF(MedicalData[Quarter],"Q3")'
$ws.Range("E104").Formula = '=COUNTIF(MedicalData[Quarter],"Q4")'
$ws.Range("F104").Formula = "=SUM(B104:E104)"

$ws.Range("A105").Value = "Ambulance Calls"
$ws.Range("B105").Formula = '=COUNTIFS(MedicalData[Quarter],"Q1",MedicalData[Hospital Transportation],"Yes")'
$ws.Range("C105").Formula = '=COUNTIFS(MedicalData[Quarter],"Q2",MedicalData[Hospital Transportation],"Yes")'
$ws.Range("D105").Formula = '=COUNTIFS(MedicalData[Quarter],"Q3",MedicalData[Hospital Transportation],"Yes")'
$ws.Range("E105").Formula = '=COUNTIFS(MedicalData[Quarter],"Q4",MedicalData[Hospital Transportation],"Yes")'
$ws.Range("F105").Formula = "=SUM(B105:E105)"

$ws.Range("A106").Value = "Average Age"
$ws.Range("B106").Formula = '=IFERROR(ROUND(AVERAGEIF(MedicalData[Quarter],"Q1",MedicalData[Age]),0),"-")'
$ws.Range("C106").Formula = '=IFERROR(ROUND(AVERAGEIF(MedicalData[Quarter],"Q2",MedicalData[Age]),0),"-")'
$ws.Range("D106").Formula = '=IFERROR(ROUND(AVERAGEIF(MedicalData[Quarter],"Q3",MedicalData[Age]),0),"-")'
$ws.Range("E106").Formula = '=IFERROR(ROUND(AVERAGEIF(MedicalData[Quarter],"Q4",MedicalData[Age]),0),"-")'
$ws.Range("F106").Formula = '=IFERROR(ROUND(AVERAGE(MedicalData[Age]),0),"-")'
$ws.Range("F106").Copy()
$ws.Range("F106").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
# F106 gains the same "total" highlight style the other quarterly-report
# total cells use (it did not have it at the old F43 location).
$ws.Range("F104").Copy()
$ws.Range("F106").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false
$ws.Range("F106").Formula = '=IFERROR(ROUND(AVERAGE(MedicalData[Age]),0),"-")'

$ws.Range("A107").Value = "Critical Cases (P1+P2)"
$ws.Range("B107").Formula = '=COUNTIFS(MedicalData[Quarter],"Q1",MedicalData[Priority Of Pt],"P1")+COUNTIFS(MedicalData[Quarter],"Q1",MedicalData[Priority Of Pt],"P2")'
$ws.Range("C107").Formula = '=COUNTIFS(MedicalData[Quarter],"Q2",MedicalData[Priority Of Pt],"P1")+COUNTIFS(MedicalData[Quarter],"Q2",MedicalData[Priority Of Pt],"P2")'
$ws.Range("D107").Formula = '=COUNTIFS(MedicalData[Quarter],"Q3",MedicalData[Priority Of Pt],"P1")+COUNTIFS(MedicalData[Quarter],"Q3",MedicalData[Priority Of Pt],"P2")'
$ws.Range("E107").Formula = '=COUNTIFS(MedicalData[Quarter],"Q4",MedicalData[Priority Of Pt],"P1")+COUNTIFS(MedicalData[Quarter],"Q4",MedicalData[Priority Of Pt],"P2")'
$ws.Range("F107").Formula = "=SUM(B107:E107)"

$ws.Range("A109").Value = "Current Quarter Highlight:"
$ws.Range("B109").Formula = '=B101&" has "&IF(B101="Q1",B104,IF(B101="Q2",C104,IF(B101="Q3",D104,E104)))&" incidents"'

# ---------------------------------------------------------------------
# 3) Clear the old Weekly/Quarterly report locations (rows 25-46) -- the
#    charts positioned at rows 10/27/44 used to be obscured by this text.
# ---------------------------------------------------------------------

$ws.Range("A25:F34").ClearContents()
$ws.Range("A37:F46").ClearContents()

# ---------------------------------------------------------------------
# 4) Lists sheet: drop the trailing blank rows 26-29 (dimension A1:E29 -> A1:E25)
# ---------------------------------------------------------------------

$wsLists = $wb.Worksheets.Item("Lists")
$wsLists.Range("A26:E29").ClearContents()
